$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.9693716918425304
$ws.Range("J2").Value = 0.9693716918425304
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.083576666666666
$ws.Range("N2").Value = 9.250729999999999
$ws.Range("O2").Value = 0.2272509363535097
$ws.Range("P2").Value = 0.2272509363535097
$ws.Range("Q2").Value = 28.23618613788222
$ws.Range("R2").Value = 254.12567524094
$ws.Range("S2").Value = 0.2202906246458008
$ws.Range("T2").Value = 0.2202906246458008

# Row 3
$ws.Range("I3").Value = 0.9693716918425304
$ws.Range("J3").Value = 0.9693716918425304
$ws.Range("M3").Value = 6.453984666666667
$ws.Range("O3").Value = 0.4756405360586227
$ws.Range("P3").Value = 0.4756405360586227
$ws.Range("Q3").Value = 59.09887513062357
$ws.Range("S3").Value = 0.4610724711480352
$ws.Range("T3").Value = 0.4610724711480351

# Row 4
$ws.Range("I4").Value = 0.9693716918425304
$ws.Range("J4").Value = 0.9693716918425304
$ws.Range("M4").Value = 4.031477000000001
$ws.Range("N4").Value = 12.094431
$ws.Range("O4").Value = 0.2971085275878677
$ws.Range("P4").Value = 0.2971085275878677
$ws.Range("Q4").Value = 36.91607094226868
$ws.Range("R4").Value = 332.2446384804181
$ws.Range("S4").Value = 0.2880085960486944
$ws.Range("T4").Value = 0.2880085960486944

# Row 5
$ws.Range("G5").Value = 0.2893236666666667
$ws.Range("H5").Value = 0.867971
$ws.Range("I5").Value = 0.03062830815746963
$ws.Range("J5").Value = 0.03062830815746962
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.083576666666666
$ws.Range("N5").Value = 9.250729999999999
$ws.Range("O5").Value = 0.2272509363535097
$ws.Range("P5").Value = 0.2272509363535097
$ws.Range("Q5").Value = 0.8921517076477777
$ws.Range("R5").Value = 8.02936536883
$ws.Range("S5").Value = 0.006960311707708811
$ws.Range("T5").Value = 0.00696031170770881

# Row 6
$ws.Range("G6").Value = 0.2893236666666667
$ws.Range("H6").Value = 0.867971
$ws.Range("I6").Value = 0.03062830815746963
$ws.Range("J6").Value = 0.03062830815746962
$ws.Range("M6").Value = 6.453984666666667
$ws.Range("O6").Value = 0.4756405360586227
$ws.Range("P6").Value = 0.4756405360586227
$ws.Range("Q6").Value = 1.867290508370445
$ws.Range("R6").Value = 16.805614575334
$ws.Range("S6").Value = 0.01456806491058754
$ws.Range("T6").Value = 0.01456806491058754

# Row 7
$ws.Range("G7").Value = 0.2893236666666667
$ws.Range("H7").Value = 0.867971
$ws.Range("I7").Value = 0.03062830815746963
$ws.Range("J7").Value = 0.03062830815746962
$ws.Range("M7").Value = 4.031477000000001
$ws.Range("N7").Value = 12.094431
$ws.Range("O7").Value = 0.2971085275878677
$ws.Range("P7").Value = 0.2971085275878677
$ws.Range("Q7").Value = 1.166401707722334
$ws.Range("R7").Value = 10.497615369501
$ws.Range("S7").Value = 0.009099931539173278
$ws.Range("T7").Value = 0.009099931539173276
